$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header columns (Wins, Losses, Ties) after the existing
# "Unnamed: 28" column (AC), copying the header formatting (style) from the
# neighboring header cell so the new headers match the existing look.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
$ws.Range("AD2:AD61").Value = 74
$ws.Range("AE2:AE61").Value = 88
$ws.Range("AF2:AF61").Value = 0
